# Commit: "deleted rare name slide"
#
# The deck contained a slide whose sole shape is a picture with
# alt-text/description "infrequent_name.png" (an example of a rare /
# infrequent baby name). That slide - along with its notes page, which
# PowerPoint removes automatically when the slide is deleted - is no
# longer part of the deck.

$p = $ppt.ActivePresentation

$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $isRareNameSlide = $false
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.Name -eq "Shape 123") {
            $isRareNameSlide = $true
        }
        if ($shape.AlternativeText -like "*infrequent_name*") {
            $isRareNameSlide = $true
        }
    }
    if ($isRareNameSlide) {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    # Fallback: this was the 11th slide in the original deck.
    $targetIndex = 11
}

$p.Slides.Item($targetIndex).Delete()
